# Basic_Excel_2cells.xlsx — "initial test seems to work to populate cell by range"
#
# Summary of the change being applied:
#  - The old single defined name "Cell_1" (-> Sheet1!$B$3) is replaced by two
#    new named ranges: "cell1rangename" (-> Sheet1!$C$3) and
#    "cell2rangename" (-> Sheet1!$F$6).
#  - The two label/value pairs that used to live at A3 and C5:D5 move to a
#    new B3:C3 / E6:F6 layout, and the label text + value text change.
#  - The two (now-unused) built-in Hyperlink / Followed Hyperlink cell
#    styles are removed from the style table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Defined names: drop the old one, add the two new range names.
# ---------------------------------------------------------------------
$wb.Names.Item("Cell_1").Delete()
$wb.Names.Add("cell1rangename", "=Sheet1!`$C`$3")
$wb.Names.Add("cell2rangename", "=Sheet1!`$F`$6")

# ---------------------------------------------------------------------
# 2. Move the two label/value pairs to their new homes.
#    Clear the old cells first so the sheet's used range shrinks back
#    down instead of leaving stale entries behind.
# ---------------------------------------------------------------------
$ws.Range("A3").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("B3").Value = "Cell1"
$ws.Range("C3").Value = "oldvalue"
$ws.Range("E6").Value = "Cell2"
$ws.Range("F6").Value = "oldvalue2"

# ---------------------------------------------------------------------
# 3. Move the active selection to the first populated cell of the new
#    layout (matches the saved <selection activeCell="C3" sqref="C3"/>).
# ---------------------------------------------------------------------
$ws.Range("C3").Select()

# Reposition the workbook window (cosmetic; matches the updated
# xWindow/yWindow saved in the workbook view).
$excel.ActiveWindow.Left = 4860
$excel.ActiveWindow.Top = 2780

# ---------------------------------------------------------------------
# 4. Drop the now-unused built-in hyperlink cell styles.
# ---------------------------------------------------------------------
$wb.Styles.Item("Followed Hyperlink").Delete()
$wb.Styles.Item("Hyperlink").Delete()
